# Fruta / hortaliza, semanal
# Re-shuffle the weekly price rows (2-13) of the active sheet.
# Row 7 is unchanged; all other rows receive the Fecha (D), Volumen (M),
# Precio minimo/maximo/promedio (N/O/P), Origen (R) and Precio $/Kg (S)
# values that used to belong to a different row, per the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2  <- old row 11
$ws.Range("D2").Value = 44616
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 3200
$ws.Range("O2").Value = 3200
$ws.Range("P2").Value = 3200
$ws.Range("R2").Value = "Región de La Araucanía"
$ws.Range("S2").Value = 3200

# Row 3  <- old row 4
$ws.Range("D3").Value = 44214
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 1800
$ws.Range("O3").Value = 1800
$ws.Range("P3").Value = 1800
$ws.Range("R3").Value = "Región de La Araucanía"
$ws.Range("S3").Value = 1800

# Row 4  <- old row 8
$ws.Range("D4").Value = 44175
$ws.Range("M4").Value = 40
$ws.Range("N4").Value = 5000
$ws.Range("O4").Value = 5000
$ws.Range("P4").Value = 5000
$ws.Range("R4").Value = "Provincia de Curicó"
$ws.Range("S4").Value = 5000

# Row 5  <- old row 13
$ws.Range("D5").Value = 44592
$ws.Range("M5").Value = 5
$ws.Range("N5").Value = 7500
$ws.Range("O5").Value = 7500
$ws.Range("P5").Value = 7500
$ws.Range("R5").Value = "Región de La Araucanía"
$ws.Range("S5").Value = 7500

# Row 6  <- old row 12
$ws.Range("D6").Value = 44323
$ws.Range("M6").Value = 20
$ws.Range("N6").Value = 3200
$ws.Range("O6").Value = 3200
$ws.Range("P6").Value = 3200
$ws.Range("R6").Value = "Región de La Araucanía"
$ws.Range("S6").Value = 3200

# Row 7 is unchanged.

# Row 8  <- old row 9
$ws.Range("D8").Value = 44176
$ws.Range("M8").Value = 20
$ws.Range("N8").Value = 3000
$ws.Range("O8").Value = 3000
$ws.Range("P8").Value = 3000
$ws.Range("R8").Value = "Región de O'Higgins"
$ws.Range("S8").Value = 3000

# Row 9  <- old row 2
$ws.Range("D9").Value = 44574
$ws.Range("M9").Value = 200
$ws.Range("N9").Value = 3000
$ws.Range("O9").Value = 3000
$ws.Range("P9").Value = 3000
$ws.Range("R9").Value = "Región de La Araucanía"
$ws.Range("S9").Value = 3000

# Row 10 <- old row 5
$ws.Range("D10").Value = 44215
$ws.Range("M10").Value = 65
$ws.Range("N10").Value = 2800
$ws.Range("O10").Value = 2800
$ws.Range("P10").Value = 2800
$ws.Range("R10").Value = "Región de La Araucanía"
$ws.Range("S10").Value = 2800

# Row 11 <- old row 10
$ws.Range("D11").Value = 44551
$ws.Range("M11").Value = 120
$ws.Range("N11").Value = 4500
$ws.Range("O11").Value = 4500
$ws.Range("P11").Value = 4500
$ws.Range("R11").Value = "Región de O'Higgins"
$ws.Range("S11").Value = 4500

# Row 12 <- old row 6
$ws.Range("D12").Value = 44567
$ws.Range("M12").Value = 80
$ws.Range("N12").Value = 2400
$ws.Range("O12").Value = 2400
$ws.Range("P12").Value = 2400
$ws.Range("R12").Value = "Región de La Araucanía"
$ws.Range("S12").Value = 2400

# Row 13 <- old row 3
$ws.Range("D13").Value = 44999
$ws.Range("M13").Value = 25
$ws.Range("N13").Value = 2500
$ws.Range("O13").Value = 2500
$ws.Range("P13").Value = 2500
$ws.Range("R13").Value = "Región de La Araucanía"
$ws.Range("S13").Value = 2500
